# Lhb-Lhcgr.xlsx was regenerated from an updated TPM matrix. For the
# Lhb -> Lhcgr pair there used to be two target-cluster rows (ECs, MuSCs);
# with the refreshed data only the MuSCs row survives, carrying new edge
# statistics, so the ECs row is dropped and row 2's numbers are refreshed
# to the new MuSCs-vs-FAPs figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now represents the (refreshed) MuSCs target-cluster figures.
$ws.Range("D2").Value = "MuSCs"

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.141084
$ws.Range("N2").Value = 6.423252
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.7192158086079999
$ws.Range("R2").Value = 6.472942277472
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# The old row 3 (ECs target-cluster data) no longer exists in the refreshed export.
$ws.Rows("3:3").Delete()
